$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <= original row 4
$ws.Range("B2").Value = "PME,met,place,patriotisme,économique"
$ws.Range("C2").Value = "un patriotisme économique"
$ws.Range("D2").Value = "un protectionnisme intelligent"
$ws.Range("E2").Value = "il dit à"
$ws.Range("F2").Value = "dit,constructeurs,américains,voulez"

# Row 3 <= original row 2
$ws.Range("B3").Value = "patriotisme,économique,jamais,mis,œuvre"
$ws.Range("C3").Value = "mis en œuvre"
$ws.Range("D3").Value = "le protectionnisme intelligent"
$ws.Range("E3").Value = "la défiscalisation de"
$ws.Range("F3").Value = "défiscalisation,heures,supplémentaires,suppression"

# Row 4 <= original row 6
$ws.Range("B4").Value = "`$,`$,`$,mettre,œuvre"
$ws.Range("C4").Value = "en œuvre de"
$ws.Range("D4").Value = "le protectionnisme intelligent"
$ws.Range("E4").Value = "à mettre en"
$ws.Range("F4").Value = "mettre,patriotisme,économique,donner"

# Row 5 <= original row 3
$ws.Range("B5").Value = "mettre,œuvre,protectionnisme,intelligent,mettre"
$ws.Range("C5").Value = "en avant de"
$ws.Range("D5").Value = "le patriotisme économique"
$ws.Range("E5").Value = "pour donner un"
$ws.Range("F5").Value = "donner,avantage,entreprises,françaises"

# Row 6 <= original row 5
$ws.Range("B6").Value = "exclusivement,TPE,PME,met,place"
$ws.Range("C6").Value = "met en place"
$ws.Range("D6").Value = "un patriotisme économique"
$ws.Range("E6").Value = "un protectionnisme intelligent"
$ws.Range("F6").Value = "protectionnisme,intelligent,dit,constructeurs"

# Row 10 <= original row 13
$ws.Range("B10").Value = "patriotisme,économique,protectionnisme,intelligent,dit"
$ws.Range("C10").Value = "il dit à"
$ws.Range("D10").Value = "les constructeurs américains"
$ws.Range("E10").Value = "si vous voulez"
$ws.Range("F10").Value = "voulez,aller,faire,voitures"

# Row 12 <= original row 15
$ws.Range("B12").Value = "construire,voiture,étranger,paierez,taxe"
$ws.Range("C12").Value = "une taxe en"
$ws.Range("D12").Value = "les réimportant"
$ws.Range("E12").Value = "à les Etats-Unis"
$ws.Range("F12").Value = "Etats-Unis,autant,évidemment,Trump"

# Row 13 <= original row 10
$ws.Range("B13").Value = "supplémentaires,suppression,travail,détaché,baisse"
$ws.Range("C13").Value = "la baisse de"
$ws.Range("D13").Value = "les charges"
$ws.Range("E13").Value = "mais exclusivement pour"
$ws.Range("F13").Value = "exclusivement,TPE,PME,met"

# Row 15 <= original row 12
$ws.Range("B15").Value = "travail,détaché,baisse,charges,exclusivement"
$ws.Range("C15").Value = "mais exclusivement pour"
$ws.Range("D15").Value = "les TPE PME ."
$ws.Range("E15").Value = "Il met en"
$ws.Range("F15").Value = "met,place,patriotisme,économique"

# Row 30 <= original row 31
$ws.Range("B30").Value = "place,politique,appelle,vœux,notamment"
$ws.Range("C30").Value = "longtemps et notamment"
$ws.Range("D30").Value = "la politique"
$ws.Range("E30").Value = "de patriotisme économique"
$ws.Range("F30").Value = "patriotisme,économique,protectionnisme,intelligent"

# Row 31 <= original row 30
$ws.Range("B31").Value = "intéresse,cause,puisqu,met,place"
$ws.Range("C31").Value = "met en place"
$ws.Range("D31").Value = "la politique"
$ws.Range("E31").Value = "que j’ appelle"
$ws.Range("F31").Value = "appelle,vœux,notamment,politique"
